$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the value to be stored as text (matches existing inlineStr/text cells)
    # without leaving a residual direct-formatted style on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "245.32"
Set-TextValue $ws.Range("D3") "23.72"
Set-TextValue $ws.Range("D4") "5.370"
Set-TextValue $ws.Range("D5") "0.05878"
Set-TextValue $ws.Range("D6") "6.477"
Set-TextValue $ws.Range("D7") "3.352"
Set-TextValue $ws.Range("D8") "0.8114"
Set-TextValue $ws.Range("D9") "0.9196"
Set-TextValue $ws.Range("D10") "0.1419"
Set-TextValue $ws.Range("D11") "0.07417"
Set-TextValue $ws.Range("D12") "0.03112"
Set-TextValue $ws.Range("D13") "0.03055"
Set-TextValue $ws.Range("D14") "0.09352"
Set-TextValue $ws.Range("D15") "3.861"
Set-TextValue $ws.Range("D16") "0.001560"
Set-TextValue $ws.Range("D17") "0.04710"
Set-TextValue $ws.Range("D18") "0.0006023"
Set-TextValue $ws.Range("D19") "0.005874"
Set-TextValue $ws.Range("D20") "0.001248"
Set-TextValue $ws.Range("D21") "0.004705"
Set-TextValue $ws.Range("D22") "0.00008812"
Set-TextValue $ws.Range("D23") "3.599"
Set-TextValue $ws.Range("D25") "0.3229"
Set-TextValue $ws.Range("D40") "0.03864"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D41") "0.1068"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D42") "0.002754"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D43") "0.003075"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
Set-TextValue $ws.Range("D44") "0.008079"
Set-TextValue $ws.Range("D45") "0.00005248"
Set-TextValue $ws.Range("D47") "0.6804"
Set-TextValue $ws.Range("D48") "0.001694"
$ws.Range("E48").Value = "47BOLOBOLO"
